# Add team record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) for the three new columns, mirroring the formatting
# of the existing header cells (bold, thin border, centered/top aligned).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerSrc = $ws.Range("AC1")
$headerSrc.Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$lastRow = 68

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 73
    $ws.Cells.Item($r, 31).Value = 89
    $ws.Cells.Item($r, 32).Value = 0
}
